# "excel loading in progress" — bump the loading-date cell by one day and
# leave the selection parked on it (A2), matching where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 held the date serial 43882 (2020-02-21); advance it to 43883 (2020-02-22).
$ws.Range("A2").Value = 43883

# Move/save the active selection to A2 (was J4).
[void]$ws.Range("A2").Select()
